# Deploy the implementation guide.
#
# Applies the changes described by the commit:
#  - rename the "Include ..." sheet to "Include #0"
#  - refresh the Metadata sheet: new Date, new Contact, insert a
#    "Jurisdiction" property (shifting Description/Purpose/Copyright/
#    Immutable down by one row) and re-append Immutable's value at
#    the new bottom row.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item(1)
$wsInclude = $wb.Worksheets.Item(2)

# 1. Rename the second sheet.
$wsInclude.Name = "Include #0"

# 2. Update the Date value (row 8).
$wsMeta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 3. Update the Contact value (row 10).
$wsMeta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 4. Insert a "Jurisdiction" property before "Description", pushing the
#    remaining rows (Description, Purpose, Copyright, Immutable) down by
#    one. Do this bottom-up so we don't clobber values we still need.

# New row 15 = old row 14 (Immutable / BooleanType[null]). Clone the
# formatting of row 14 (the last existing data row) onto the brand-new
# row 15 before writing its values, so it picks up the same cell style.
$wsMeta.Range("A14:B14").Copy() | Out-Null
$wsMeta.Range("A15:B15").PasteSpecial(-4122) | Out-Null
$wsMeta.Range("A15").Value = "Immutable"
$wsMeta.Range("B15").Value = "BooleanType[null]"

# Row 14 becomes Copyright (old row 13's label), value stays blank.
$wsMeta.Range("A14").Value = "Copyright"
$wsMeta.Range("B14").ClearContents()

# Row 13 becomes Purpose (old row 12's label), value stays blank.
$wsMeta.Range("A13").Value = "Purpose"

# Row 12 becomes Description (old row 11's label/value).
$wsMeta.Range("A12").Value = "Description"
$wsMeta.Range("B12").Value = "Population ValueSet"

# Row 11 becomes the new Jurisdiction property with an empty value.
$wsMeta.Range("A11").Value = "Jurisdiction"
$wsMeta.Range("B11").Value = ""
